$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dev Log")

# --- Insert a new development-log entry as row 4 (pushes existing rows down) ---
$ws.Rows("4:4").Insert()

# Copy the formatting (row height / number formats / fonts / fills) from what is
# now row 5 (the former row 4) so the freshly inserted row matches its siblings.
$ws.Range("A5:G5").Copy()
$ws.Range("A4:G4").PasteSpecial(-4122)
$ws.Rows(4).RowHeight = 297.95
$ws.Application.CutCopyMode = $false

# --- Populate the new log entry ---
# (cell order matters for the shared-string table layout: PROGRESS first,
#  then FIX/STATUS, then PERCEPTION, matching the author's entry order)
$ws.Range("B4").Value = 45363
$ws.Range("C4").Value = 0.56666666666666665
$ws.Range("F4").Value = "Begin logic for user to enter tiles"
$ws.Range("D4").Value = "Work on logic for user to enter tiles"
$ws.Range("E4").Value = "Onwards and upwards"
$ws.Range("G4").Value = 0.54

# --- Extend the conditional formatting so it still covers the full log table ---
$exprRule = $ws.Range("B4:F12").FormatConditions.Item(1)
$exprRule.ModifyAppliesToRange($ws.Range("B4:F13"))

$barRule = $ws.Range("G4:G12").FormatConditions.Item(1)
$barRule.ModifyAppliesToRange($ws.Range("G4:G13"))

# --- Reset zoom on the Dev Log sheet to 50% ---
$ws.Activate()
$ws.Range("F4").Select()
$excel.ActiveWindow.Zoom = 50
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
